{"js": "// Update the header date and the 25 division problems in the table, in\n// document order, matching the author's commit (output generated at c986bee).\n\nconst dateMap = [\"2024-11-04 Monday\", \"2024-11-05 Tuesday\"];\n\n// Sequential (row-major) replacements for every non-blank table cell.\nconst cellReplacements = [\n  \"36\u00f79=\", \"19\u00f76=\", \"51\u00f72=\", \"58\u00f79=\", \"67\u00f75=\",\n  \"62\u00f77=\", \"46\u00f72=\", \"17\u00f74=\", \"12\u00f76=\", \"51\u00f79=\",\n  \"25\u00f76=\", \"99\u00f74=\", \"68\u00f72=\", \"11\u00f76=\", \"48\u00f76=\",\n  \"54\u00f78=\", \"59\u00f78=\", \"65\u00f79=\", \"30\u00f77=\", \"92\u00f76=\",\n  \"17\u00f75=\", \"88\u00f79=\", \"70\u00f77=\", \"57\u00f75=\", \"87\u00f77=\",\n];\n\n// 1. Update the first paragraph (the date line).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(dateMap[1], Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Update every populated cell in the table, left-to-right / top-to-bottom.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst colCount = table.values[0].length;\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const current = table.values[r][c];\n    if (current === \"\") {\n      continue; // blank spacer cell, nothing to replace\n    }\n    const newText = cellReplacements[i];\n    i++;\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    const cellParagraph = cell.body.paragraphs.items[0];\n    cellParagraph.getRange().insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Update the header date and the 25 division problems in the table, in\n# document order, matching the author's commit (output generated at c986bee).\n\n$d = $word.ActiveDocument\n\n# 1. Update the first paragraph (the date line).\n$d.Paragraphs.Item(1).Range.Text = \"2024-11-05 Tuesday\"\n\n# 2. Update every populated cell in the table, left-to-right / top-to-bottom.\n#    The table has 20 rows of 5 columns; only rows 1, 5, 9, 13, 17 carry\n#    the division problems, the rest are blank spacer rows.\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$replacements = @(\n    @(\"36\u00f79=\", \"19\u00f76=\", \"51\u00f72=\", \"58\u00f79=\", \"67\u00f75=\"),\n    @(\"62\u00f77=\", \"46\u00f72=\", \"17\u00f74=\", \"12\u00f76=\", \"51\u00f79=\"),\n    @(\"25\u00f76=\", \"99\u00f74=\", \"68\u00f72=\", \"11\u00f76=\", \"48\u00f76=\"),\n    @(\"54\u00f78=\", \"59\u00f78=\", \"65\u00f79=\", \"30\u00f77=\", \"92\u00f76=\"),\n    @(\"17\u00f75=\", \"88\u00f79=\", \"70\u00f77=\", \"57\u00f75=\", \"87\u00f77=\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $row = $dataRows[$i]\n    $rowValues = $replacements[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($row, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
